# Commit: "updating error message in the siq"
#
# Appends three new SIQ question/answer rows (Q58-Q60) about validation
# error messages to the bottom of the requirements table on Sheet2 (rows
# 79-81), matching the layout/format of the existing rows directly above
# (75-78).
#
# Columns in this table (see header rows 67-78):
#   B = ID              E = Status            H = Who asked
#   C = Question         F = Who answered      I = Date asked
#   D = Answer           G = Who answered      J = Due date

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 79
$lastNewRow  = 81
$srcRow      = $firstNewRow - 1   # row 78: last existing row, used as the formatting template

# Clone formatting (fonts/fills/borders/alignment/row height) from the last
# existing data row onto the new rows before filling in content, so the
# appended rows look identical to the rest of the table.
$ws.Range("B$srcRow`:J$srcRow").Copy()
$ws.Range("B$firstNewRow`:J$lastNewRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- fill content column by column (mirrors how the rows were authored) --

$ws.Range("B79").Value = "BANK_SYS_SIQ_Q58"
$ws.Range("B80").Value = "BANK_SYS_SIQ_Q59"
$ws.Range("B81").Value = "BANK_SYS_SIQ_Q60"

$ws.Range("C79").Value = "what is the error message shall be displayed if the customer leave a mandatory field empty ?"
$ws.Range("C80").Value = "what is the error message shall be displayed if the customer doesn't achieve length constraint ?"
$ws.Range("C81").Value = "what is the error message shall be displayed if the customer doesn't achieve the using chararcters constraints  ?"

$ws.Range("D79").Value = 'error message beyound the field that " this field is mandatory "'
$ws.Range("D80").Value = 'error message beyound the field that "length must be between x and y "'
$ws.Range("D81").Value = 'error message beyound the field that " invalid data format "'

$ws.Range("E$firstNewRow`:E$lastNewRow").Value = "Agreed with your proposal."
$ws.Range("F$firstNewRow`:F$lastNewRow").Value = "answered"
$ws.Range("G$firstNewRow`:G$lastNewRow").Value = "khadija mostafa"
$ws.Range("H$firstNewRow`:H$lastNewRow").Value = "Marwan"
$ws.Range("I$firstNewRow`:J$lastNewRow").Value = '"18/5/2019"'

# Match the row heights of the rest of the table.
$ws.Rows("$firstNewRow`:$lastNewRow").RowHeight = $ws.Rows("$srcRow").RowHeight

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("J84").Select()
